# Workbook currently has:
#   Worksheets.Item(1) -> physical sheet1.xml (rId1), name "hotel_info",
#                          holds the hotel_info header (A1:I1) + 1 data row (A2:I2)
#   Worksheets.Item(2) -> physical sheet2.xml (rId2), name "review_info",
#                          holds the review_info header (A1:Y1), no data rows
#
# Target:
#   rId1 (Worksheets.Item(1)) -> name "review_info", holding the review_info
#                                 header only (A1:Y1), no data rows
#   rId2 (Worksheets.Item(2)) -> name "hotel_info", holding the hotel_info
#                                 header + 1 data row, with a new "State"
#                                 column inserted right after "Hotel_Name"
#                                 (before "City"), value "Louisiana"

$wb = $excel.ActiveWorkbook
$sheetA = $wb.Worksheets.Item(1)
$sheetB = $wb.Worksheets.Item(2)

# --- capture the current hotel_info content (lives on $sheetA) before wiping ---
$hotelA2 = $sheetA.Range("A2").Value()   # 52917 (number)
$hotelB2 = $sheetA.Range("B2").Value()   # Suburban Extended Stay Hotel Avondale
$hotelC2 = $sheetA.Range("C2").Value()   # Avondale
$hotelD2 = $sheetA.Range("D2").Value()   # 70094 (number)
$hotelE2 = $sheetA.Range("E2").Value()   # TA review URL
$hotelF2 = $sheetA.Range("F2").Value()   # Tripadvisor hotel name
$hotelG2 = $sheetA.Range("G2").Value()   # "72" (English_Reviews_num, stored as text)
$hotelH2 = $sheetA.Range("H2").Value()   # "2"  (Local_Rank, stored as text)

# --- the review_info header, currently sitting on $sheetB ---
$reviewHeader = @("STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL","Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title","review_content","review_rating","trip_month","trip_purpose","value","rooms","Location","Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text")

# --- wipe both sheets clean before rebuilding ---
$sheetA.Cells.Clear()
$sheetB.Cells.Clear()

# --- rename tabs to match target (go through temp names to dodge the
#     name-collision check, since the final names are each other's
#     current name) ---
$sheetA.Name = "TEMP_SHEET_A"
$sheetB.Name = "TEMP_SHEET_B"
$sheetA.Name = "review_info"
$sheetB.Name = "hotel_info"

# --- rebuild sheetA ($sheetA / rId1) as review_info: header row only ---
$reviewCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")
for ($i = 0; $i -lt $reviewCols.Length; $i++) {
    $cell = $reviewCols[$i] + "1"
    $sheetA.Range($cell).Value = $reviewHeader[$i]
}

# --- rebuild sheetB ($sheetB / rId2) as hotel_info: header + data row, with new "State" column ---
$sheetB.Range("A1").Value = "STR"
$sheetB.Range("B1").Value = "Hotel_Name"
$sheetB.Range("C1").Value = "State"
$sheetB.Range("D1").Value = "City"
$sheetB.Range("E1").Value = "Zip"
$sheetB.Range("F1").Value = "TA_ReviewURL"
$sheetB.Range("G1").Value = "Tripadvisor_Hotel_Name"
$sheetB.Range("H1").Value = "English_Reviews_num"
$sheetB.Range("I1").Value = "Local_Rank"
$sheetB.Range("J1").Value = "Total_Reviews_num"

$sheetB.Range("A2").Value = $hotelA2
$sheetB.Range("B2").Value = $hotelB2
$sheetB.Range("C2").Value = "Louisiana"
$sheetB.Range("D2").Value = $hotelC2
$sheetB.Range("E2").Value = $hotelD2
$sheetB.Range("F2").Value = $hotelE2
$sheetB.Range("G2").Value = $hotelF2

# English_Reviews_num / Local_Rank / Total_Reviews_num are text cells even
# though they look numeric ("72", "2") in the source data - force text
# storage via NumberFormat "@" so Excel doesn't silently convert them to
# numbers, then reset the style so no stray number-format sticks around.
$textCells = $sheetB.Range("H2:J2")
$textCells.NumberFormat = "@"
$sheetB.Range("H2").Value = $hotelG2
$sheetB.Range("I2").Value = $hotelH2
$sheetB.Range("J2").Value = $hotelG2
$textCells.Style = "Normal"

Write-Output "done"
